# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows (new date 2022-11-24 / serial 44889) above the
# existing historical rows, pushing the old rows 41-46 down to 44-49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 41-46 down by 3 rows (rows 41-43 become 44-46, etc.)
$ws.Rows("41:43").Insert()

# --- Fill the 3 new rows (41, 42, 43) with this week's data ---

$ws.Range("A41:A43").Value = 1
$ws.Range("B41:B43").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C41:C43").Value = "Arica y Parinacota"

$ws.Range("D41:D43").Value = 44889
$ws.Range("D41:D43").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("E41:E43").Value = 15
$ws.Range("F41:F43").Value = 100112028
$ws.Range("G41:G43").Value = "Sandia"
$ws.Range("H41:H43").Value = "Sin especificar"

$ws.Range("I41").Value = "Primera"
$ws.Range("I42").Value = "Segunda"
$ws.Range("I43").Value = "Tercera"

$ws.Range("J41").Value = 540
$ws.Range("K41").Value = 480
$ws.Range("L41").Value = 500
$ws.Range("M41").Value = 489
$ws.Range("P41").Value = 489

$ws.Range("J42").Value = 350
$ws.Range("K42").Value = 480
$ws.Range("L42").Value = 500
$ws.Range("M42").Value = 491
$ws.Range("P42").Value = 491

$ws.Range("J43").Value = 190
$ws.Range("K43").Value = 480
$ws.Range("L43").Value = 500
$ws.Range("M43").Value = 491
$ws.Range("P43").Value = 491

$ws.Range("N41:N43").Value = '$/kilo (volumen en unidades)'
$ws.Range("O41:O43").Value = "Perú"
$ws.Range("Q41:Q43").Value = 1
$ws.Range("R41:R43").Value = "Hortaliza"
